# RPA datasets push 2023-11-10
# Rebuild data rows 2-32 to reflect the refreshed IPO underwriting dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCols = @("B", "F", "G")
foreach ($col in $dateCols) {
  $ws.Range($col + "2:" + $col + "32").NumberFormat = "@"
}

# Row 2
$ws.Range("A2").Value = 'CS'
$ws.Range("B2").Value = '2023-09-21'
$ws.Range("C2").Value = '두산로보틱스'
$ws.Range("D2").Value = '한국, 미래'
$ws.Range("E2").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F2").Value = '2023-09-26'
$ws.Range("G2").Value = '2023-10-05'
$ws.Range("H2").Value = 42120
$ws.Range("I2").Value = 16200000
$ws.Range("J2").Value = 26000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 10

# Row 3
$ws.Range("A3").Value = 'IBK'
$ws.Range("B3").Value = '2023-10-31'
$ws.Range("C3").Value = '비아이매트릭스'
$ws.Range("D3").Value = 'IBK'
$ws.Range("E3").Value = 'IBK'
$ws.Range("F3").Value = '2023-11-03'
$ws.Range("G3").Value = '2023-11-09'
$ws.Range("H3").Value = 15600
$ws.Range("I3").Value = 1200000
$ws.Range("J3").Value = 13000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100

# Row 4
$ws.Range("A4").Value = 'KB'
$ws.Range("B4").Value = '2023-10-26'
$ws.Range("C4").Value = '쏘닉스'
$ws.Range("D4").Value = 'KB'
$ws.Range("E4").Value = 'KB'
$ws.Range("F4").Value = '2023-10-31'
$ws.Range("G4").Value = '2023-11-07'
$ws.Range("H4").Value = 27000
$ws.Range("I4").Value = 3600000
$ws.Range("J4").Value = 7500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100

# Row 5
$ws.Range("A5").Value = 'KB'
$ws.Range("B5").Value = '2023-10-24'
$ws.Range("C5").Value = 'KB제27호스팩'
$ws.Range("D5").Value = 'KB'
$ws.Range("E5").Value = 'KB'
$ws.Range("F5").Value = '2023-10-27'
$ws.Range("G5").Value = '2023-11-03'
$ws.Range("H5").Value = 25000
$ws.Range("I5").Value = 12500000
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100

# Row 6
$ws.Range("A6").Value = 'KB'
$ws.Range("B6").Value = '2023-09-21'
$ws.Range("C6").Value = '두산로보틱스'
$ws.Range("D6").Value = '한국, 미래'
$ws.Range("E6").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F6").Value = '2023-09-26'
$ws.Range("G6").Value = '2023-10-05'
$ws.Range("H6").Value = 42120
$ws.Range("I6").Value = 16200000
$ws.Range("J6").Value = 26000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 10

# Row 7
$ws.Range("A7").Value = 'KB'
$ws.Range("B7").Value = '2023-09-19'
$ws.Range("C7").Value = '한싹'
$ws.Range("D7").Value = 'KB'
$ws.Range("E7").Value = 'KB'
$ws.Range("F7").Value = '2023-09-22'
$ws.Range("G7").Value = '2023-10-04'
$ws.Range("H7").Value = 18750
$ws.Range("I7").Value = 1500000
$ws.Range("J7").Value = 12500
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 100

# Row 8
$ws.Range("A8").Value = 'NH'
$ws.Range("B8").Value = '2023-09-21'
$ws.Range("C8").Value = '두산로보틱스'
$ws.Range("D8").Value = '한국, 미래'
$ws.Range("E8").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F8").Value = '2023-09-26'
$ws.Range("G8").Value = '2023-10-05'
$ws.Range("H8").Value = 42120
$ws.Range("I8").Value = 16200000
$ws.Range("J8").Value = 26000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 10

# Row 9
$ws.Range("A9").Value = 'NH'
$ws.Range("B9").Value = '2023-10-31'
$ws.Range("C9").Value = '메가터치'
$ws.Range("D9").Value = 'NH'
$ws.Range("E9").Value = 'NH'
$ws.Range("F9").Value = '2023-11-03'
$ws.Range("G9").Value = '2023-11-10'
$ws.Range("H9").Value = 24960
$ws.Range("I9").Value = 5200000
$ws.Range("J9").Value = 4800
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100

# Row 10
$ws.Range("A10").Value = 'NH'
$ws.Range("B10").Value = '2023-10-23'
$ws.Range("C10").Value = '유진테크놀로지'
$ws.Range("D10").Value = 'NH'
$ws.Range("E10").Value = 'NH'
$ws.Range("F10").Value = '2023-10-26'
$ws.Range("G10").Value = '2023-11-02'
$ws.Range("H10").Value = 17841.194
$ws.Range("I10").Value = 1049482
$ws.Range("J10").Value = 17000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 100

# Row 11
$ws.Range("A11").Value = '대신'
$ws.Range("B11").Value = '2023-10-31'
$ws.Range("C11").Value = '컨텍'
$ws.Range("D11").Value = '대신'
$ws.Range("E11").Value = '대신'
$ws.Range("F11").Value = '2023-11-03'
$ws.Range("G11").Value = '2023-11-09'
$ws.Range("H11").Value = 46350
$ws.Range("I11").Value = 2060000
$ws.Range("J11").Value = 22500
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 100

# Row 12
$ws.Range("A12").Value = '미래'
$ws.Range("B12").Value = '2023-09-18'
$ws.Range("C12").Value = '밀리의서재'
$ws.Range("D12").Value = '미래'
$ws.Range("E12").Value = '미래'
$ws.Range("F12").Value = '2023-09-21'
$ws.Range("G12").Value = '2023-09-27'
$ws.Range("H12").Value = 34500
$ws.Range("I12").Value = 1500000
$ws.Range("J12").Value = 23000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 100

# Row 13
$ws.Range("A13").Value = '미래'
$ws.Range("B13").Value = '2023-10-05'
$ws.Range("C13").Value = '퓨릿'
$ws.Range("D13").Value = '미래'
$ws.Range("E13").Value = '미래'
$ws.Range("F13").Value = '2023-10-11'
$ws.Range("G13").Value = '2023-10-18'
$ws.Range("H13").Value = 44265.9
$ws.Range("I13").Value = 4137000
$ws.Range("J13").Value = 10700
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 100

# Row 14
$ws.Range("A14").Value = '미래'
$ws.Range("B14").Value = '2023-10-31'
$ws.Range("C14").Value = '큐로셀'
$ws.Range("D14").Value = '미래, 삼성'
$ws.Range("E14").Value = '미래, 삼성'
$ws.Range("F14").Value = '2023-11-03'
$ws.Range("G14").Value = '2023-11-09'
$ws.Range("H14").Value = 16000
$ws.Range("I14").Value = 1600000
$ws.Range("J14").Value = 20000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 50

# Row 15
$ws.Range("A15").Value = '미래'
$ws.Range("B15").Value = '2023-09-21'
$ws.Range("C15").Value = '두산로보틱스'
$ws.Range("D15").Value = '한국, 미래'
$ws.Range("E15").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F15").Value = '2023-09-26'
$ws.Range("G15").Value = '2023-10-05'
$ws.Range("H15").Value = 126360
$ws.Range("I15").Value = 16200000
$ws.Range("J15").Value = 26000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 30

# Row 16
$ws.Range("A16").Value = '미래'
$ws.Range("B16").Value = '2023-10-10'
$ws.Range("C16").Value = '신성에스티'
$ws.Range("D16").Value = '미래'
$ws.Range("E16").Value = '미래'
$ws.Range("F16").Value = '2023-10-13'
$ws.Range("G16").Value = '2023-10-19'
$ws.Range("H16").Value = 52000
$ws.Range("I16").Value = 2000000
$ws.Range("J16").Value = 26000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 100

# Row 17
$ws.Range("A17").Value = '삼성'
$ws.Range("B17").Value = '2023-10-31'
$ws.Range("C17").Value = '큐로셀'
$ws.Range("D17").Value = '미래, 삼성'
$ws.Range("E17").Value = '미래, 삼성'
$ws.Range("F17").Value = '2023-11-03'
$ws.Range("G17").Value = '2023-11-09'
$ws.Range("H17").Value = 16000
$ws.Range("I17").Value = 1600000
$ws.Range("J17").Value = 20000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 50

# Row 18
$ws.Range("A18").Value = '삼성'
$ws.Range("B18").Value = '2023-09-19'
$ws.Range("C18").Value = '레뷰코퍼레이션'
$ws.Range("D18").Value = '삼성'
$ws.Range("E18").Value = '삼성'
$ws.Range("F18").Value = '2023-09-22'
$ws.Range("G18").Value = '2023-10-06'
$ws.Range("H18").Value = 33600
$ws.Range("I18").Value = 2240000
$ws.Range("J18").Value = 15000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 100

# Row 19
$ws.Range("A19").Value = '신영'
$ws.Range("B19").Value = '2023-09-21'
$ws.Range("C19").Value = '두산로보틱스'
$ws.Range("D19").Value = '한국, 미래'
$ws.Range("E19").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F19").Value = '2023-09-26'
$ws.Range("G19").Value = '2023-10-05'
$ws.Range("H19").Value = 12636
$ws.Range("I19").Value = 16200000
$ws.Range("J19").Value = 26000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 3

# Row 20
$ws.Range("A20").Value = '신영'
$ws.Range("B20").Value = '2023-09-14'
$ws.Range("C20").Value = '인스웨이브시스템즈'
$ws.Range("D20").Value = '신영'
$ws.Range("E20").Value = '신영'
$ws.Range("F20").Value = '2023-09-19'
$ws.Range("G20").Value = '2023-09-25'
$ws.Range("H20").Value = 26400
$ws.Range("I20").Value = 1100000
$ws.Range("J20").Value = 24000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 100

# Row 21
$ws.Range("A21").Value = '신한'
$ws.Range("B21").Value = '2023-09-19'
$ws.Range("C21").Value = '신한제11호스팩'
$ws.Range("D21").Value = '신한'
$ws.Range("E21").Value = '신한'
$ws.Range("F21").Value = '2023-09-22'
$ws.Range("G21").Value = '2023-10-04'
$ws.Range("H21").Value = 36000
$ws.Range("I21").Value = 18000000
$ws.Range("J21").Value = 2000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 100

# Row 22
$ws.Range("A22").Value = '신한'
$ws.Range("B22").Value = '2023-10-23'
$ws.Range("C22").Value = '유투바이오'
$ws.Range("D22").Value = '신한'
$ws.Range("E22").Value = '신한'
$ws.Range("F22").Value = '2023-10-26'
$ws.Range("G22").Value = '2023-11-02'
$ws.Range("H22").Value = 4966.368
$ws.Range("I22").Value = 1128720
$ws.Range("J22").Value = 4400
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 100

# Row 23
$ws.Range("A23").Value = '유비에스'
$ws.Range("B23").Value = '2023-09-21'
$ws.Range("C23").Value = '두산로보틱스'
$ws.Range("D23").Value = '한국, 미래'
$ws.Range("E23").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F23").Value = '2023-09-26'
$ws.Range("G23").Value = '2023-10-05'
$ws.Range("H23").Value = 4212
$ws.Range("I23").Value = 16200000
$ws.Range("J23").Value = 26000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1

# Row 24
$ws.Range("A24").Value = '유안타'
$ws.Range("B24").Value = '2023-09-18'
$ws.Range("C24").Value = '아이엠티'
$ws.Range("D24").Value = '유안타'
$ws.Range("E24").Value = '유안타, 유진'
$ws.Range("F24").Value = '2023-09-21'
$ws.Range("G24").Value = '2023-10-10'
$ws.Range("H24").Value = 15484
$ws.Range("I24").Value = 1580000
$ws.Range("J24").Value = 14000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 70

# Row 25
$ws.Range("A25").Value = '유진'
$ws.Range("B25").Value = '2023-09-18'
$ws.Range("C25").Value = '아이엠티'
$ws.Range("D25").Value = '유안타'
$ws.Range("E25").Value = '유안타, 유진'
$ws.Range("F25").Value = '2023-09-21'
$ws.Range("G25").Value = '2023-10-10'
$ws.Range("H25").Value = 6636
$ws.Range("I25").Value = 1580000
$ws.Range("J25").Value = 14000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 30

# Row 26
$ws.Range("A26").Value = '키움'
$ws.Range("B26").Value = '2023-09-21'
$ws.Range("C26").Value = '두산로보틱스'
$ws.Range("D26").Value = '한국, 미래'
$ws.Range("E26").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F26").Value = '2023-09-26'
$ws.Range("G26").Value = '2023-10-05'
$ws.Range("H26").Value = 12636
$ws.Range("I26").Value = 16200000
$ws.Range("J26").Value = 26000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 3

# Row 27
$ws.Range("A27").Value = '키움'
$ws.Range("B27").Value = '2023-10-16'
$ws.Range("C27").Value = '워트'
$ws.Range("D27").Value = '키움'
$ws.Range("E27").Value = '키움'
$ws.Range("F27").Value = '2023-10-19'
$ws.Range("G27").Value = '2023-10-26'
$ws.Range("H27").Value = 26000
$ws.Range("I27").Value = 4000000
$ws.Range("J27").Value = 6500
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 100

# Row 28
$ws.Range("A28").Value = '하나'
$ws.Range("B28").Value = '2023-10-10'
$ws.Range("C28").Value = '에스엘에스바이오'
$ws.Range("D28").Value = '하나'
$ws.Range("E28").Value = '하나'
$ws.Range("F28").Value = '2023-10-13'
$ws.Range("G28").Value = '2023-10-20'
$ws.Range("H28").Value = 5390
$ws.Range("I28").Value = 770000
$ws.Range("J28").Value = 7000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 100

# Row 29
$ws.Range("A29").Value = '하나'
$ws.Range("B29").Value = '2023-09-21'
$ws.Range("C29").Value = '두산로보틱스'
$ws.Range("D29").Value = '한국, 미래'
$ws.Range("E29").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F29").Value = '2023-09-26'
$ws.Range("G29").Value = '2023-10-05'
$ws.Range("H29").Value = 12636
$ws.Range("I29").Value = 16200000
$ws.Range("J29").Value = 26000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 3

# Row 30
$ws.Range("A30").Value = '한국'
$ws.Range("B30").Value = '2023-10-18'
$ws.Range("C30").Value = '퀄리타스반도체'
$ws.Range("D30").Value = '한국'
$ws.Range("E30").Value = '한국'
$ws.Range("F30").Value = '2023-10-23'
$ws.Range("G30").Value = '2023-10-27'
$ws.Range("H30").Value = 30600
$ws.Range("I30").Value = 1800000
$ws.Range("J30").Value = 17000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 100

# Row 31
$ws.Range("A31").Value = '한국'
$ws.Range("B31").Value = '2023-09-21'
$ws.Range("C31").Value = '두산로보틱스'
$ws.Range("D31").Value = '한국, 미래'
$ws.Range("E31").Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Range("F31").Value = '2023-09-26'
$ws.Range("G31").Value = '2023-10-05'
$ws.Range("H31").Value = 126360
$ws.Range("I31").Value = 16200000
$ws.Range("J31").Value = 26000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 30

# Row 32
$ws.Range("A32").Value = '현대차'
$ws.Range("B32").Value = '2023-09-25'
$ws.Range("C32").Value = '에이치엠씨제6호스팩'
$ws.Range("D32").Value = '현대차'
$ws.Range("E32").Value = '현대차'
$ws.Range("F32").Value = '2023-10-04'
$ws.Range("G32").Value = '2023-10-13'
$ws.Range("H32").Value = 8000
$ws.Range("I32").Value = 4000000
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 100
